$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 0  # H40: 1156.5 -> 0
$ws.Cells.Item(40, 9).Value = 0  # I40: 1000 -> 0
$ws.Cells.Item(40, 10).Value = 0  # J40: 1313 -> 0
$ws.Cells.Item(40, 11).Value = 0  # K40: 1000 -> 0
$ws.Cells.Item(40, 12).Value = 0  # L40: 1313 -> 0
$ws.Cells.Item(40, 13).Value = $null  # M40: -825 -> (removed)
$ws.Cells.Item(40, 14).Value = $null  # N40: -1663 -> (removed)
$ws.Cells.Item(41, 8).Value = 469.0625  # H41: 454.2 -> 469.0625
$ws.Cells.Item(41, 9).Value = 470.33334  # I41: 405.8889 -> 470.33334
$ws.Cells.Item(41, 10).Value = 467.42856  # J41: 526.6667 -> 467.42856
$ws.Cells.Item(41, 11).Value = 470.33334  # K41: 405.8889 -> 470.33334
$ws.Cells.Item(41, 12).Value = 467.42856  # L41: 526.6667 -> 467.42856
$ws.Cells.Item(41, 13).Value = -30.33334000000002  # M41: 34.11110000000002 -> -30.33334000000002
$ws.Cells.Item(41, 14).Value = -1347.42856  # N41: -1406.6667 -> -1347.42856
$ws.Cells.Item(80, 8).Value = 1865.6  # H80: 1377.7142 -> 1865.6
$ws.Cells.Item(80, 9).Value = 500  # I80: 343.2 -> 500
$ws.Cells.Item(80, 10).Value = 2450.8572  # J80: 1952.4445 -> 2450.8572
$ws.Cells.Item(80, 11).Value = 1500  # K80: 1029.6 -> 1500
$ws.Cells.Item(80, 12).Value = 7352.571599999999  # L80: 5857.333500000001 -> 7352.571599999999
$ws.Cells.Item(80, 13).Value = -502  # M80: -31.59999999999991 -> -502
$ws.Cells.Item(80, 14).Value = -9348.571599999999  # N80: -7853.333500000001 -> -9348.571599999999
$ws.Cells.Item(83, 8).Value = 1865.6  # H83: 1377.7142 -> 1865.6
$ws.Cells.Item(83, 9).Value = 500  # I83: 343.2 -> 500
$ws.Cells.Item(83, 10).Value = 2450.8572  # J83: 1952.4445 -> 2450.8572
$ws.Cells.Item(83, 11).Value = 4500  # K83: 3088.8 -> 4500
$ws.Cells.Item(83, 12).Value = 22057.7148  # L83: 17572.0005 -> 22057.7148
$ws.Cells.Item(83, 13).Value = 492  # M83: 1903.2 -> 492
$ws.Cells.Item(83, 14).Value = -32041.7148  # N83: -27556.0005 -> -32041.7148
$ws.Cells.Item(98, 8).Value = 1087.5  # H98: 1008.8461 -> 1087.5
$ws.Cells.Item(98, 9).Value = 1497.5  # I98: 1292.8572 -> 1497.5
$ws.Cells.Item(98, 11).Value = 1497.5  # K98: 1292.8572 -> 1497.5
$ws.Cells.Item(98, 13).Value = 0.5  # M98: 205.1428000000001 -> 0.5
$ws.Cells.Item(111, 8).Value = 1614.9231  # H111: 1640.9286 -> 1614.9231
$ws.Cells.Item(111, 10).Value = 2178  # J111: 2153.125 -> 2178
$ws.Cells.Item(111, 12).Value = 6534  # L111: 6459.375 -> 6534
$ws.Cells.Item(111, 14).Value = -12668  # N111: -12593.375 -> -12668
$ws.Cells.Item(122, 8).Value = 1087.5  # H122: 1008.8461 -> 1087.5
$ws.Cells.Item(122, 9).Value = 1497.5  # I122: 1292.8572 -> 1497.5
$ws.Cells.Item(122, 11).Value = 4492.5  # K122: 3878.5716 -> 4492.5
$ws.Cells.Item(122, 13).Value = -2042.5  # M122: -1428.5716 -> -2042.5
$ws.Cells.Item(138, 8).Value = 2377.1765  # H138: 2323 -> 2377.1765
$ws.Cells.Item(138, 10).Value = 2203.3333  # J138: 1683 -> 2203.3333
$ws.Cells.Item(138, 12).Value = 6609.999899999999  # L138: 5049 -> 6609.999899999999
$ws.Cells.Item(138, 14).Value = -16889.9999  # N138: -15329 -> -16889.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5194.3335  # H32: 6008.8887 -> 5194.3335
$ws.Cells.Item(32, 9).Value = 3380.625  # I32: 4089.923 -> 3380.625
$ws.Cells.Item(32, 11).Value = 3380.625  # K32: 4089.923 -> 3380.625
$ws.Cells.Item(32, 13).Value = -3093.625  # M32: -3802.923 -> -3093.625
$ws.Cells.Item(122, 8).Value = 3376.647  # H122: 3425.1875 -> 3376.647
$ws.Cells.Item(122, 9).Value = 3123.2856  # I122: 3163.5386 -> 3123.2856
$ws.Cells.Item(122, 11).Value = 9369.856800000001  # K122: 9490.6158 -> 9369.856800000001
$ws.Cells.Item(122, 13).Value = -6919.856800000001  # M122: -7040.6158 -> -6919.856800000001
$ws.Cells.Item(132, 8).Value = 0  # H132: 2000 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 1000 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 3000 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 3000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 9000 -> 0
$ws.Cells.Item(132, 13).Value = $null  # M132: -470 -> (removed)
$ws.Cells.Item(132, 14).Value = $null  # N132: -14060 -> (removed)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2428.6072  # H86: 2515.6296 -> 2428.6072
$ws.Cells.Item(86, 9).Value = 1732.8422  # I86: 1824.7222 -> 1732.8422
$ws.Cells.Item(86, 11).Value = 1732.8422  # K86: 1824.7222 -> 1732.8422
$ws.Cells.Item(86, 13).Value = -609.8422  # M86: -701.7221999999999 -> -609.8422
$ws.Cells.Item(89, 8).Value = 2428.6072  # H89: 2515.6296 -> 2428.6072
$ws.Cells.Item(89, 9).Value = 1732.8422  # I89: 1824.7222 -> 1732.8422
$ws.Cells.Item(89, 11).Value = 8664.210999999999  # K89: 9123.610999999999 -> 8664.210999999999
$ws.Cells.Item(89, 13).Value = -3048.210999999999  # M89: -3507.610999999999 -> -3048.210999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 1495.5  # H7: 1572.9565 -> 1495.5
$ws.Cells.Item(7, 9).Value = 632.05554  # I7: 698.1875 -> 632.05554
$ws.Cells.Item(7, 10).Value = 3438.25  # J7: 3572.4285 -> 3438.25
$ws.Cells.Item(7, 11).Value = 632.05554  # K7: 698.1875 -> 632.05554
$ws.Cells.Item(7, 12).Value = 3438.25  # L7: 3572.4285 -> 3438.25
$ws.Cells.Item(7, 13).Value = -519.05554  # M7: -585.1875 -> -519.05554
$ws.Cells.Item(7, 14).Value = -3664.25  # N7: -3798.4285 -> -3664.25
$ws.Cells.Item(140, 8).Value = 100780  # H140: 0 -> 100780
$ws.Cells.Item(140, 10).Value = 100780  # J140: 0 -> 100780
$ws.Cells.Item(140, 12).Value = 100780  # L140: 0 -> 100780
$ws.Cells.Item(140, 14).Value = -111140  # N140: None -> -111140

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 654.0909  # H122: 524.6875 -> 654.0909
$ws.Cells.Item(122, 9).Value = 459.66666  # I122: 447.42856 -> 459.66666
$ws.Cells.Item(122, 10).Value = 887.4  # J122: 584.7778 -> 887.4
$ws.Cells.Item(122, 11).Value = 4136.99994  # K122: 4026.85704 -> 4136.99994
$ws.Cells.Item(122, 12).Value = 7986.599999999999  # L122: 5263.000199999999 -> 7986.599999999999
$ws.Cells.Item(122, 13).Value = -1686.99994  # M122: -1576.85704 -> -1686.99994
$ws.Cells.Item(122, 14).Value = -12886.6  # N122: -10163.0002 -> -12886.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 0  # H70: 4639.5 -> 0
$ws.Cells.Item(70, 9).Value = 0  # I70: 4639.5 -> 0
$ws.Cells.Item(70, 11).Value = 0  # K70: 4639.5 -> 0
$ws.Cells.Item(70, 13).Value = $null  # M70: -4369.5 -> (removed)
$ws.Cells.Item(73, 8).Value = 0  # H73: 4639.5 -> 0
$ws.Cells.Item(73, 9).Value = 0  # I73: 4639.5 -> 0
$ws.Cells.Item(73, 11).Value = 0  # K73: 4639.5 -> 0
$ws.Cells.Item(73, 13).Value = $null  # M73: -3703.5 -> (removed)
$ws.Cells.Item(80, 8).Value = 2476.6667  # H80: 2708.8 -> 2476.6667
$ws.Cells.Item(80, 9).Value = 2199.6667  # I80: 2399.5 -> 2199.6667
$ws.Cells.Item(80, 10).Value = 2753.6667  # J80: 2915 -> 2753.6667
$ws.Cells.Item(80, 11).Value = 2199.6667  # K80: 2399.5 -> 2199.6667
$ws.Cells.Item(80, 12).Value = 2753.6667  # L80: 2915 -> 2753.6667
$ws.Cells.Item(80, 13).Value = -1201.6667  # M80: -1401.5 -> -1201.6667
$ws.Cells.Item(80, 14).Value = -4749.6667  # N80: -4911 -> -4749.6667
$ws.Cells.Item(83, 8).Value = 2476.6667  # H83: 2708.8 -> 2476.6667
$ws.Cells.Item(83, 9).Value = 2199.6667  # I83: 2399.5 -> 2199.6667
$ws.Cells.Item(83, 10).Value = 2753.6667  # J83: 2915 -> 2753.6667
$ws.Cells.Item(83, 11).Value = 10998.3335  # K83: 11997.5 -> 10998.3335
$ws.Cells.Item(83, 12).Value = 13768.3335  # L83: 14575 -> 13768.3335
$ws.Cells.Item(83, 13).Value = -6006.333500000001  # M83: -7005.5 -> -6006.333500000001
$ws.Cells.Item(83, 14).Value = -23752.3335  # N83: -24559 -> -23752.3335
$ws.Cells.Item(126, 8).Value = 5397.8  # H126: 5597.5 -> 5397.8
$ws.Cells.Item(126, 10).Value = 5699.6665  # J126: 6250 -> 5699.6665
$ws.Cells.Item(126, 12).Value = 17098.9995  # L126: 18750 -> 17098.9995
$ws.Cells.Item(126, 14).Value = -22038.9995  # N126: -23690 -> -22038.9995
$ws.Cells.Item(132, 8).Value = 4113.3335  # H132: 1761.0667 -> 4113.3335
$ws.Cells.Item(132, 9).Value = 4470  # I132: 1644 -> 4470
$ws.Cells.Item(132, 11).Value = 13410  # K132: 4932 -> 13410
$ws.Cells.Item(132, 13).Value = -10880  # M132: -2402 -> -10880

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 0  # H22: 2949.75 -> 0
$ws.Cells.Item(22, 9).Value = 0  # I22: 300 -> 0
$ws.Cells.Item(22, 10).Value = 0  # J22: 3833 -> 0
$ws.Cells.Item(22, 11).Value = 0  # K22: 300 -> 0
$ws.Cells.Item(22, 12).Value = 0  # L22: 3833 -> 0
$ws.Cells.Item(22, 13).Value = $null  # M22: -5 -> (removed)
$ws.Cells.Item(22, 14).Value = $null  # N22: -4423 -> (removed)
$ws.Cells.Item(27, 8).Value = 0  # H27: 2949.75 -> 0
$ws.Cells.Item(27, 9).Value = 0  # I27: 300 -> 0
$ws.Cells.Item(27, 10).Value = 0  # J27: 3833 -> 0
$ws.Cells.Item(27, 11).Value = 0  # K27: 300 -> 0
$ws.Cells.Item(27, 12).Value = 0  # L27: 3833 -> 0
$ws.Cells.Item(27, 13).Value = $null  # M27: -193 -> (removed)
$ws.Cells.Item(27, 14).Value = $null  # N27: -4047 -> (removed)
$ws.Cells.Item(38, 8).Value = 31000  # H38: 0 -> 31000
$ws.Cells.Item(38, 10).Value = 31000  # J38: 0 -> 31000
$ws.Cells.Item(38, 12).Value = 31000  # L38: 0 -> 31000
$ws.Cells.Item(38, 14).Value = -31820  # N38: None -> -31820
$ws.Cells.Item(46, 8).Value = 1584.04  # H46: 1678.1305 -> 1584.04
$ws.Cells.Item(46, 10).Value = 1989.3572  # J46: 2237.25 -> 1989.3572
$ws.Cells.Item(46, 12).Value = 1989.3572  # L46: 2237.25 -> 1989.3572
$ws.Cells.Item(46, 14).Value = -2365.3572  # N46: -2613.25 -> -2365.3572
$ws.Cells.Item(61, 8).Value = 933  # H61: 900 -> 933
$ws.Cells.Item(61, 9).Value = 899.5  # I61: 900 -> 899.5
$ws.Cells.Item(61, 10).Value = 1000  # J61: 0 -> 1000
$ws.Cells.Item(61, 11).Value = 899.5  # K61: 900 -> 899.5
$ws.Cells.Item(61, 12).Value = 1000  # L61: 0 -> 1000
$ws.Cells.Item(61, 13).Value = -697.5  # M61: -698 -> -697.5
$ws.Cells.Item(61, 14).Value = -1404  # N61: None -> -1404
$ws.Cells.Item(82, 8).Value = 1184.7778  # H82: 1093.9 -> 1184.7778
$ws.Cells.Item(82, 9).Value = 931  # I82: 802.7143 -> 931
$ws.Cells.Item(82, 10).Value = 1387.8  # J82: 1773.3334 -> 1387.8
$ws.Cells.Item(82, 11).Value = 931  # K82: 802.7143 -> 931
$ws.Cells.Item(82, 12).Value = 1387.8  # L82: 1773.3334 -> 1387.8
$ws.Cells.Item(82, 13).Value = -570  # M82: -441.7143 -> -570
$ws.Cells.Item(82, 14).Value = -2109.8  # N82: -2495.3334 -> -2109.8
$ws.Cells.Item(85, 8).Value = 1184.7778  # H85: 1093.9 -> 1184.7778
$ws.Cells.Item(85, 9).Value = 931  # I85: 802.7143 -> 931
$ws.Cells.Item(85, 10).Value = 1387.8  # J85: 1773.3334 -> 1387.8
$ws.Cells.Item(85, 11).Value = 931  # K85: 802.7143 -> 931
$ws.Cells.Item(85, 12).Value = 1387.8  # L85: 1773.3334 -> 1387.8
$ws.Cells.Item(85, 13).Value = 317  # M85: 445.2857 -> 317
$ws.Cells.Item(85, 14).Value = -3883.8  # N85: -4269.3334 -> -3883.8
$ws.Cells.Item(100, 8).Value = 4297.8887  # H100: 4424.75 -> 4297.8887
$ws.Cells.Item(100, 9).Value = 3613  # I100: 3668 -> 3613
$ws.Cells.Item(100, 11).Value = 3613  # K100: 3668 -> 3613
$ws.Cells.Item(100, 13).Value = -3072  # M100: -3127 -> -3072
$ws.Cells.Item(113, 8).Value = 933  # H113: 900 -> 933
$ws.Cells.Item(113, 9).Value = 899.5  # I113: 900 -> 899.5
$ws.Cells.Item(113, 10).Value = 1000  # J113: 0 -> 1000
$ws.Cells.Item(113, 11).Value = 899.5  # K113: 900 -> 899.5
$ws.Cells.Item(113, 12).Value = 1000  # L113: 0 -> 1000
$ws.Cells.Item(113, 13).Value = 1270.5  # M113: 1270 -> 1270.5
$ws.Cells.Item(113, 14).Value = -5340  # N113: None -> -5340
$ws.Cells.Item(132, 8).Value = 2804.4075  # H132: 2820 -> 2804.4075
$ws.Cells.Item(132, 9).Value = 2783.2727  # I132: 2801.5715 -> 2783.2727
$ws.Cells.Item(132, 11).Value = 8349.8181  # K132: 8404.7145 -> 8349.8181
$ws.Cells.Item(132, 13).Value = -5819.8181  # M132: -5874.7145 -> -5819.8181
$ws.Cells.Item(136, 8).Value = 4268.5625  # H136: 4440.1333 -> 4268.5625
$ws.Cells.Item(136, 10).Value = 5385.5713  # J136: 6000.6665 -> 5385.5713
$ws.Cells.Item(136, 12).Value = 16156.7139  # L136: 18001.9995 -> 16156.7139
$ws.Cells.Item(136, 14).Value = -21256.7139  # N136: -23101.9995 -> -21256.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4299.5  # H62: 4332.1665 -> 4299.5
$ws.Cells.Item(62, 9).Value = 3799.25  # I62: 4066 -> 3799.25
$ws.Cells.Item(62, 10).Value = 5300  # J62: 4598.3335 -> 5300
$ws.Cells.Item(62, 11).Value = 3799.25  # K62: 4066 -> 3799.25
$ws.Cells.Item(62, 12).Value = 5300  # L62: 4598.3335 -> 5300
$ws.Cells.Item(62, 13).Value = -3175.25  # M62: -3442 -> -3175.25
$ws.Cells.Item(62, 14).Value = -6548  # N62: -5846.3335 -> -6548
$ws.Cells.Item(65, 8).Value = 4299.5  # H65: 4332.1665 -> 4299.5
$ws.Cells.Item(65, 9).Value = 3799.25  # I65: 4066 -> 3799.25
$ws.Cells.Item(65, 10).Value = 5300  # J65: 4598.3335 -> 5300
$ws.Cells.Item(65, 11).Value = 18996.25  # K65: 20330 -> 18996.25
$ws.Cells.Item(65, 12).Value = 26500  # L65: 22991.6675 -> 26500
$ws.Cells.Item(65, 13).Value = -15876.25  # M65: -17210 -> -15876.25
$ws.Cells.Item(65, 14).Value = -32740  # N65: -29231.6675 -> -32740
